# Auto-generated script: replace 100 arithmetic expressions in the table cells.
# Uses MatchWholeWord=$true to avoid accidental substring matches (e.g. "1+31=" inside "41+31=").
$d = $word.ActiveDocument

$d.Content.Find.Execute("33+10=", $false, $true, $false, $false, $false, $true, 1, $false, "36-7=", 2) | Out-Null
$d.Content.Find.Execute("10+64=", $false, $true, $false, $false, $false, $true, 1, $false, "39-23=", 2) | Out-Null
$d.Content.Find.Execute("37+27=", $false, $true, $false, $false, $false, $true, 1, $false, "58-50=", 2) | Out-Null
$d.Content.Find.Execute("16+9=", $false, $true, $false, $false, $false, $true, 1, $false, "45-44=", 2) | Out-Null
$d.Content.Find.Execute("61-15=", $false, $true, $false, $false, $false, $true, 1, $false, "40-12=", 2) | Out-Null
$d.Content.Find.Execute("39-14=", $false, $true, $false, $false, $false, $true, 1, $false, "68+28=", 2) | Out-Null
$d.Content.Find.Execute("39+38=", $false, $true, $false, $false, $false, $true, 1, $false, "36+57=", 2) | Out-Null
$d.Content.Find.Execute("56+1=", $false, $true, $false, $false, $false, $true, 1, $false, "60-8=", 2) | Out-Null
$d.Content.Find.Execute("56+19=", $false, $true, $false, $false, $false, $true, 1, $false, "88-65=", 2) | Out-Null
$d.Content.Find.Execute("1+31=", $false, $true, $false, $false, $false, $true, 1, $false, "81-79=", 2) | Out-Null
$d.Content.Find.Execute("88-41=", $false, $true, $false, $false, $false, $true, 1, $false, "96-94=", 2) | Out-Null
$d.Content.Find.Execute("76-12=", $false, $true, $false, $false, $false, $true, 1, $false, "14-8=", 2) | Out-Null
$d.Content.Find.Execute("85-82=", $false, $true, $false, $false, $false, $true, 1, $false, "18+2=", 2) | Out-Null
$d.Content.Find.Execute("93-58=", $false, $true, $false, $false, $false, $true, 1, $false, "10+81=", 2) | Out-Null
$d.Content.Find.Execute("7+81=", $false, $true, $false, $false, $false, $true, 1, $false, "10+87=", 2) | Out-Null
$d.Content.Find.Execute("30+11=", $false, $true, $false, $false, $false, $true, 1, $false, "64-5=", 2) | Out-Null
$d.Content.Find.Execute("37-30=", $false, $true, $false, $false, $false, $true, 1, $false, "22+26=", 2) | Out-Null
$d.Content.Find.Execute("87-33=", $false, $true, $false, $false, $false, $true, 1, $false, "20+45=", 2) | Out-Null
$d.Content.Find.Execute("1+79=", $false, $true, $false, $false, $false, $true, 1, $false, "30-2=", 2) | Out-Null
$d.Content.Find.Execute("12+75=", $false, $true, $false, $false, $false, $true, 1, $false, "4+73=", 2) | Out-Null
$d.Content.Find.Execute("96-92=", $false, $true, $false, $false, $false, $true, 1, $false, "40-20=", 2) | Out-Null
$d.Content.Find.Execute("50-16=", $false, $true, $false, $false, $false, $true, 1, $false, "27-27=", 2) | Out-Null
$d.Content.Find.Execute("98-90=", $false, $true, $false, $false, $false, $true, 1, $false, "63-52=", 2) | Out-Null
$d.Content.Find.Execute("81-71=", $false, $true, $false, $false, $false, $true, 1, $false, "75+17=", 2) | Out-Null
$d.Content.Find.Execute("2+91=", $false, $true, $false, $false, $false, $true, 1, $false, "42+8=", 2) | Out-Null
$d.Content.Find.Execute("47-15=", $false, $true, $false, $false, $false, $true, 1, $false, "14+9=", 2) | Out-Null
$d.Content.Find.Execute("51-22=", $false, $true, $false, $false, $false, $true, 1, $false, "49-18=", 2) | Out-Null
$d.Content.Find.Execute("46+46=", $false, $true, $false, $false, $false, $true, 1, $false, "0+85=", 2) | Out-Null
$d.Content.Find.Execute("67-5=", $false, $true, $false, $false, $false, $true, 1, $false, "18+6=", 2) | Out-Null
$d.Content.Find.Execute("95-31=", $false, $true, $false, $false, $false, $true, 1, $false, "84-35=", 2) | Out-Null
$d.Content.Find.Execute("89-62=", $false, $true, $false, $false, $false, $true, 1, $false, "37+13=", 2) | Out-Null
$d.Content.Find.Execute("93-60=", $false, $true, $false, $false, $false, $true, 1, $false, "33+33=", 2) | Out-Null
$d.Content.Find.Execute("69+30=", $false, $true, $false, $false, $false, $true, 1, $false, "94-29=", 2) | Out-Null
$d.Content.Find.Execute("57-36=", $false, $true, $false, $false, $false, $true, 1, $false, "19-9=", 2) | Out-Null
$d.Content.Find.Execute("52-17=", $false, $true, $false, $false, $false, $true, 1, $false, "52+23=", 2) | Out-Null
$d.Content.Find.Execute("45-16=", $false, $true, $false, $false, $false, $true, 1, $false, "35+0=", 2) | Out-Null
$d.Content.Find.Execute("11+16=", $false, $true, $false, $false, $false, $true, 1, $false, "78-49=", 2) | Out-Null
$d.Content.Find.Execute("59+17=", $false, $true, $false, $false, $false, $true, 1, $false, "59-6=", 2) | Out-Null
$d.Content.Find.Execute("3+69=", $false, $true, $false, $false, $false, $true, 1, $false, "96-34=", 2) | Out-Null
$d.Content.Find.Execute("62-60=", $false, $true, $false, $false, $false, $true, 1, $false, "48+50=", 2) | Out-Null
$d.Content.Find.Execute("33-5=", $false, $true, $false, $false, $false, $true, 1, $false, "89-74=", 2) | Out-Null
$d.Content.Find.Execute("39-3=", $false, $true, $false, $false, $false, $true, 1, $false, "18+49=", 2) | Out-Null
$d.Content.Find.Execute("20+53=", $false, $true, $false, $false, $false, $true, 1, $false, "38-36=", 2) | Out-Null
$d.Content.Find.Execute("42+43=", $false, $true, $false, $false, $false, $true, 1, $false, "78-44=", 2) | Out-Null
$d.Content.Find.Execute("72-60=", $false, $true, $false, $false, $false, $true, 1, $false, "46-32=", 2) | Out-Null
$d.Content.Find.Execute("92-33=", $false, $true, $false, $false, $false, $true, 1, $false, "61-35=", 2) | Out-Null
$d.Content.Find.Execute("10+25=", $false, $true, $false, $false, $false, $true, 1, $false, "38+60=", 2) | Out-Null
$d.Content.Find.Execute("77+1=", $false, $true, $false, $false, $false, $true, 1, $false, "30+20=", 2) | Out-Null
$d.Content.Find.Execute("6+78=", $false, $true, $false, $false, $false, $true, 1, $false, "37+28=", 2) | Out-Null
$d.Content.Find.Execute("72-67=", $false, $true, $false, $false, $false, $true, 1, $false, "36+8=", 2) | Out-Null
$d.Content.Find.Execute("94-25=", $false, $true, $false, $false, $false, $true, 1, $false, "95-27=", 2) | Out-Null
$d.Content.Find.Execute("61-32=", $false, $true, $false, $false, $false, $true, 1, $false, "7+24=", 2) | Out-Null
$d.Content.Find.Execute("69-27=", $false, $true, $false, $false, $false, $true, 1, $false, "16+72=", 2) | Out-Null
$d.Content.Find.Execute("57+24=", $false, $true, $false, $false, $false, $true, 1, $false, "97-45=", 2) | Out-Null
$d.Content.Find.Execute("43+36=", $false, $true, $false, $false, $false, $true, 1, $false, "27+37=", 2) | Out-Null
$d.Content.Find.Execute("76+21=", $false, $true, $false, $false, $false, $true, 1, $false, "33-3=", 2) | Out-Null
$d.Content.Find.Execute("81-8=", $false, $true, $false, $false, $false, $true, 1, $false, "89-84=", 2) | Out-Null
$d.Content.Find.Execute("98-13=", $false, $true, $false, $false, $false, $true, 1, $false, "82-0=", 2) | Out-Null
$d.Content.Find.Execute("29-27=", $false, $true, $false, $false, $false, $true, 1, $false, "52+8=", 2) | Out-Null
$d.Content.Find.Execute("60-36=", $false, $true, $false, $false, $false, $true, 1, $false, "66-61=", 2) | Out-Null
$d.Content.Find.Execute("25-16=", $false, $true, $false, $false, $false, $true, 1, $false, "83-60=", 2) | Out-Null
$d.Content.Find.Execute("75+21=", $false, $true, $false, $false, $false, $true, 1, $false, "30+60=", 2) | Out-Null
$d.Content.Find.Execute("2+93=", $false, $true, $false, $false, $false, $true, 1, $false, "0+1=", 2) | Out-Null
$d.Content.Find.Execute("41+31=", $false, $true, $false, $false, $false, $true, 1, $false, "83-16=", 2) | Out-Null
$d.Content.Find.Execute("86-83=", $false, $true, $false, $false, $false, $true, 1, $false, "6+12=", 2) | Out-Null
$d.Content.Find.Execute("1+94=", $false, $true, $false, $false, $false, $true, 1, $false, "40+48=", 2) | Out-Null
$d.Content.Find.Execute("55-24=", $false, $true, $false, $false, $false, $true, 1, $false, "60-33=", 2) | Out-Null
$d.Content.Find.Execute("75-41=", $false, $true, $false, $false, $false, $true, 1, $false, "51-34=", 2) | Out-Null
$d.Content.Find.Execute("77-9=", $false, $true, $false, $false, $false, $true, 1, $false, "68-32=", 2) | Out-Null
$d.Content.Find.Execute("11+77=", $false, $true, $false, $false, $false, $true, 1, $false, "55-8=", 2) | Out-Null
$d.Content.Find.Execute("10+24=", $false, $true, $false, $false, $false, $true, 1, $false, "54-40=", 2) | Out-Null
$d.Content.Find.Execute("62-38=", $false, $true, $false, $false, $false, $true, 1, $false, "20-14=", 2) | Out-Null
$d.Content.Find.Execute("5+16=", $false, $true, $false, $false, $false, $true, 1, $false, "35+62=", 2) | Out-Null
$d.Content.Find.Execute("13-3=", $false, $true, $false, $false, $false, $true, 1, $false, "7+75=", 2) | Out-Null
$d.Content.Find.Execute("77+20=", $false, $true, $false, $false, $false, $true, 1, $false, "7+41=", 2) | Out-Null
$d.Content.Find.Execute("28+41=", $false, $true, $false, $false, $false, $true, 1, $false, "95+4=", 2) | Out-Null
$d.Content.Find.Execute("1+61=", $false, $true, $false, $false, $false, $true, 1, $false, "21+17=", 2) | Out-Null
$d.Content.Find.Execute("29+19=", $false, $true, $false, $false, $false, $true, 1, $false, "54-24=", 2) | Out-Null
$d.Content.Find.Execute("45-27=", $false, $true, $false, $false, $false, $true, 1, $false, "45+21=", 2) | Out-Null
$d.Content.Find.Execute("81+9=", $false, $true, $false, $false, $false, $true, 1, $false, "51+16=", 2) | Out-Null
$d.Content.Find.Execute("51+45=", $false, $true, $false, $false, $false, $true, 1, $false, "7+35=", 2) | Out-Null
$d.Content.Find.Execute("30-1=", $false, $true, $false, $false, $false, $true, 1, $false, "43-15=", 2) | Out-Null
$d.Content.Find.Execute("37-27=", $false, $true, $false, $false, $false, $true, 1, $false, "91-58=", 2) | Out-Null
$d.Content.Find.Execute("38-37=", $false, $true, $false, $false, $false, $true, 1, $false, "69-55=", 2) | Out-Null
$d.Content.Find.Execute("90-21=", $false, $true, $false, $false, $false, $true, 1, $false, "87-71=", 2) | Out-Null
$d.Content.Find.Execute("82+1=", $false, $true, $false, $false, $false, $true, 1, $false, "23-20=", 2) | Out-Null
$d.Content.Find.Execute("94-54=", $false, $true, $false, $false, $false, $true, 1, $false, "78-66=", 2) | Out-Null
$d.Content.Find.Execute("12+50=", $false, $true, $false, $false, $false, $true, 1, $false, "9+37=", 2) | Out-Null
$d.Content.Find.Execute("51+24=", $false, $true, $false, $false, $false, $true, 1, $false, "56+25=", 2) | Out-Null
$d.Content.Find.Execute("74-50=", $false, $true, $false, $false, $false, $true, 1, $false, "71-41=", 2) | Out-Null
$d.Content.Find.Execute("39+22=", $false, $true, $false, $false, $false, $true, 1, $false, "0+74=", 2) | Out-Null
$d.Content.Find.Execute("49+33=", $false, $true, $false, $false, $false, $true, 1, $false, "63-47=", 2) | Out-Null
$d.Content.Find.Execute("65-61=", $false, $true, $false, $false, $false, $true, 1, $false, "93-11=", 2) | Out-Null
$d.Content.Find.Execute("31+47=", $false, $true, $false, $false, $false, $true, 1, $false, "9+61=", 2) | Out-Null
$d.Content.Find.Execute("38-8=", $false, $true, $false, $false, $false, $true, 1, $false, "44-18=", 2) | Out-Null
$d.Content.Find.Execute("21+59=", $false, $true, $false, $false, $false, $true, 1, $false, "14+8=", 2) | Out-Null
$d.Content.Find.Execute("63-16=", $false, $true, $false, $false, $false, $true, 1, $false, "32+60=", 2) | Out-Null
$d.Content.Find.Execute("27+11=", $false, $true, $false, $false, $false, $true, 1, $false, "60-2=", 2) | Out-Null
$d.Content.Find.Execute("86-11=", $false, $true, $false, $false, $false, $true, 1, $false, "3+58=", 2) | Out-Null
$d.Content.Find.Execute("85-29=", $false, $true, $false, $false, $false, $true, 1, $false, "28+46=", 2) | Out-Null
